$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target layout:
#   A: Name   B: DOB   C: Place   D: ID   E: Gender
# The existing "Age" column (B) keeps its numeric values but is relabeled
# "DOB". Two new columns, ID and Gender, are appended after Place (C),
# matched to each row by name/value.

# New header "ID" first, then relabel B1 "Age " -> "DOB", then "Gender" -
# this keeps the shared-string insertion order ID, DOB, Gender.
$ws.Range("D1").Value = "ID"
$ws.Range("B1").Value = "DOB"
$ws.Range("E1").Value = "Gender"

# Add ID values (column D)
$ws.Range("D2").Value = 23
$ws.Range("D3").Value = 34
$ws.Range("D4").Value = 35
$ws.Range("D5").Value = 65

# Add Gender values (column E), matched by column name/value per row
$ws.Range("E2").Value = "m"
$ws.Range("E3").Value = "m"
$ws.Range("E4").Value = "f"
$ws.Range("E5").Value = "f"

# Update selection to reflect the new active cell used after edits
$ws.Range("G3").Select()
